# "Fruta / hortaliza, semanal" - weekly update for the
# Fruta, Femacal de La Calera - Mango data set.
#
# A new weekly price record is inserted as row 205 (pushing the existing
# rows 205-232 down to 206-233), and the worksheet dimension grows from
# A1:T232 to A1:T233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 205; rows 205-232 shift down to 206-233
$ws.Rows.Item(205).Insert()

# Populate the new row 205 with the new weekly record
$ws.Cells.Item(205, 1).Value = 3
$ws.Cells.Item(205, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(205, 3).Value = "Coquimbo"
$ws.Cells.Item(205, 4).Value = 44491
$ws.Cells.Item(205, 5).Value = 5
$ws.Cells.Item(205, 6).Value = "Fruta"
$ws.Cells.Item(205, 7).Value = 100108
$ws.Cells.Item(205, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(205, 9).Value = 100108002
$ws.Cells.Item(205, 10).Value = "Mango"
$ws.Cells.Item(205, 11).Value = "Sin especificar"
$ws.Cells.Item(205, 12).Value = "Primera"
$ws.Cells.Item(205, 13).Value = 228
$ws.Cells.Item(205, 14).Value = 7500
$ws.Cells.Item(205, 15).Value = 7500
$ws.Cells.Item(205, 16).Value = 7500
$ws.Cells.Item(205, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(205, 18).Value = "Perú"
$ws.Cells.Item(205, 19).Value = 1875
$ws.Cells.Item(205, 20).Value = 4
